# Applies the "2018 dollars" -> "2019 dollars" refresh of the OCCF workbook
# (RMI 3.0 script data drop), per the commit "Drop in all data files from
# 3.0 RMI script".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("About")

# --- Update the underlying conversion factor and its label (About!A26/B26) ---
# Dependent formulas on the other sheets (=10^9*About!$A$26, etc.) recalc
# automatically from this new value.
$ws.Range("A26").Value = 0.89805481563188172
$ws.Range("B26").Value = "2019 dollars per 2012 dollar"

# --- Update the remaining text labels that referenced the old 2018-dollar base year ---
$ws.Range("B29").Value = 'which in this case is "2012 dollars per 2019 dollar."'
$ws.Range("A21").Value = "million 2019 dollars"
$ws.Range("A18").Value = "billion 2019 dollars"

# --- Restore the active selection to A19, matching the refreshed workbook ---
$ws.Activate()
[void]$ws.Range("A19").Select()

[void]$wb.Application.Calculate()
